# Auto-generated: apply market-price / leve-profit data refresh
# across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 313.25
$ws.Range("I33").Value = 313.25
$ws.Range("K33").Value = 313.25
$ws.Range("M33").Value = -84.25
$ws.Range("H137").Value = 418158.75
$ws.Range("I137").Value = 3437387.8
$ws.Range("J137").Value = 6445.727
$ws.Range("K137").Value = 10312163.4
$ws.Range("L137").Value = 19337.181
$ws.Range("M137").Value = -10309613.4
$ws.Range("N137").Value = -24437.181
$ws.Range("H138").Value = 169858.06
$ws.Range("J138").Value = 6313.492
$ws.Range("L138").Value = 18940.476
$ws.Range("N138").Value = -29220.476

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3909.4722
$ws.Range("I32").Value = 4276.3438
$ws.Range("J32").Value = 974.5
$ws.Range("K32").Value = 4276.3438
$ws.Range("L32").Value = 974.5
$ws.Range("M32").Value = -3989.3438
$ws.Range("N32").Value = -1548.5
$ws.Range("H61").Value = 7621.5
$ws.Range("I61").Value = 8091.8335
$ws.Range("K61").Value = 8091.8335
$ws.Range("M61").Value = -7879.8335
$ws.Range("H74").Value = 2388.0557
$ws.Range("I74").Value = 2160.5386
$ws.Range("K74").Value = 2160.5386
$ws.Range("M74").Value = -1286.5386
$ws.Range("H77").Value = 2388.0557
$ws.Range("I77").Value = 2160.5386
$ws.Range("K77").Value = 10802.693
$ws.Range("M77").Value = -6434.692999999999
$ws.Range("H88").Value = 111114160
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 111114160
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 111114160
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -111114972
$ws.Range("H91").Value = 111114160
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 111114160
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 111114160
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -111116968
$ws.Range("H102").Value = 2938.2122
$ws.Range("I102").Value = 2483.1538
$ws.Range("K102").Value = 2483.1538
$ws.Range("M102").Value = -861.1538
$ws.Range("H132").Value = 3482.2932
$ws.Range("I132").Value = 2304.725
$ws.Range("K132").Value = 6914.174999999999
$ws.Range("M132").Value = -4384.174999999999
$ws.Range("H136").Value = 7621.5
$ws.Range("I136").Value = 8091.8335
$ws.Range("K136").Value = 24275.5005
$ws.Range("M136").Value = -21725.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5746.25
$ws.Range("I20").Value = 2092.5
$ws.Range("J20").Value = 9400
$ws.Range("K20").Value = 2092.5
$ws.Range("L20").Value = 9400
$ws.Range("M20").Value = -1845.5
$ws.Range("N20").Value = -9894
$ws.Range("H86").Value = 4613.7646
$ws.Range("I86").Value = 4942.067
$ws.Range("K86").Value = 4942.067
$ws.Range("M86").Value = -3819.067
$ws.Range("H89").Value = 4613.7646
$ws.Range("I89").Value = 4942.067
$ws.Range("K89").Value = 24710.335
$ws.Range("M89").Value = -19094.335
$ws.Range("H94").Value = 1103181.4
$ws.Range("I94").Value = 1835611.1
$ws.Range("K94").Value = 1835611.1
$ws.Range("M94").Value = -1835160.1
$ws.Range("H99").Value = 18629.096
$ws.Range("I99").Value = 20929.834
$ws.Range("K99").Value = 20929.834
$ws.Range("M99").Value = -19431.834
$ws.Range("H105").Value = 7000
$ws.Range("I105").Value = 7000
$ws.Range("K105").Value = 7000
$ws.Range("M105").Value = -5253
$ws.Range("H107").Value = 3605.4285
$ws.Range("I107").Value = 3343.476
$ws.Range("J107").Value = 4391.2856
$ws.Range("K107").Value = 3343.476
$ws.Range("L107").Value = 4391.2856
$ws.Range("M107").Value = -1423.476
$ws.Range("N107").Value = -8231.285599999999
$ws.Range("H134").Value = 6839.121
$ws.Range("I134").Value = 7141.3667
$ws.Range("K134").Value = 21424.1001
$ws.Range("M134").Value = -18889.1001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2989.4792
$ws.Range("J31").Value = 3152.238
$ws.Range("L31").Value = 3152.238
$ws.Range("N31").Value = -3742.238
$ws.Range("H34").Value = 2989.4792
$ws.Range("J34").Value = 3152.238
$ws.Range("L34").Value = 3152.238
$ws.Range("N34").Value = -3556.238
$ws.Range("H132").Value = 21652.166
$ws.Range("I132").Value = 9472.6
$ws.Range("J132").Value = 82550
$ws.Range("K132").Value = 28417.8
$ws.Range("L132").Value = 247650
$ws.Range("M132").Value = -25887.8
$ws.Range("N132").Value = -252710

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 633.6
$ws.Range("I8").Value = 633.6
$ws.Range("K8").Value = 1900.8
$ws.Range("M8").Value = -1761.8
$ws.Range("H21").Value = 795.1177
$ws.Range("I21").Value = 844.6875
$ws.Range("J21").Value = 2
$ws.Range("K21").Value = 2534.0625
$ws.Range("L21").Value = 6
$ws.Range("M21").Value = -2361.0625
$ws.Range("N21").Value = -352
$ws.Range("H68").Value = 7814928.5
$ws.Range("I68").Value = 2571.4285
$ws.Range("K68").Value = 7714.2855
$ws.Range("M68").Value = -6903.2855
$ws.Range("H71").Value = 7814928.5
$ws.Range("I71").Value = 2571.4285
$ws.Range("K71").Value = 23142.8565
$ws.Range("M71").Value = -19086.8565
$ws.Range("H80").Value = 305749.75
$ws.Range("I80").Value = 2999
$ws.Range("J80").Value = 406666.66
$ws.Range("K80").Value = 8997
$ws.Range("L80").Value = 1219999.98
$ws.Range("M80").Value = -8061
$ws.Range("N80").Value = -1221871.98
$ws.Range("H83").Value = 305749.75
$ws.Range("I83").Value = 2999
$ws.Range("J83").Value = 406666.66
$ws.Range("K83").Value = 26991
$ws.Range("L83").Value = 3659999.94
$ws.Range("M83").Value = -22311
$ws.Range("N83").Value = -3669359.94
$ws.Range("H114").Value = 12061.223
$ws.Range("I114").Value = 344
$ws.Range("K114").Value = 1032
$ws.Range("M114").Value = 2222
$ws.Range("H117").Value = 16571
$ws.Range("I117").Value = 3565
$ws.Range("J117").Value = 20287
$ws.Range("K117").Value = 10695
$ws.Range("L117").Value = 60861
$ws.Range("M117").Value = -7253
$ws.Range("N117").Value = -67745
$ws.Range("H131").Value = 43482864
$ws.Range("J131").Value = 2132.111
$ws.Range("L131").Value = 6396.333
$ws.Range("N131").Value = -16476.333
$ws.Range("H132").Value = 13933725
$ws.Range("J132").Value = 16720125
$ws.Range("L132").Value = 150481125
$ws.Range("N132").Value = -150486185

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H102").Value = 9227.5
$ws.Range("I102").Value = 10316.471
$ws.Range("J102").Value = 3056.6667
$ws.Range("K102").Value = 10316.471
$ws.Range("L102").Value = 3056.6667
$ws.Range("M102").Value = -8694.471
$ws.Range("N102").Value = -6300.6667
$ws.Range("H126").Value = 16937.842
$ws.Range("J126").Value = 15356.875
$ws.Range("L126").Value = 46070.625
$ws.Range("N126").Value = -51010.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5828.5
$ws.Range("I100").Value = 6609.1816
$ws.Range("J100").Value = 2966
$ws.Range("K100").Value = 6609.1816
$ws.Range("L100").Value = 2966
$ws.Range("M100").Value = -6068.1816
$ws.Range("N100").Value = -4048
$ws.Range("H132").Value = 1330731.9
$ws.Range("I132").Value = 1624449.5
$ws.Range("J132").Value = 9002.5
$ws.Range("K132").Value = 4873348.5
$ws.Range("L132").Value = 27007.5
$ws.Range("M132").Value = -4870818.5
$ws.Range("N132").Value = -32067.5
$ws.Range("H136").Value = 6200.4814
$ws.Range("J136").Value = 8763.75
$ws.Range("L136").Value = 26291.25
$ws.Range("N136").Value = -31391.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 13770.538
$ws.Range("I45").Value = 14569
$ws.Range("J45").Value = 13704
$ws.Range("K45").Value = 14569
$ws.Range("L45").Value = 13704
$ws.Range("M45").Value = -14078
$ws.Range("N45").Value = -14686
$ws.Range("H46").Value = 84989.75
$ws.Range("J46").Value = 84989.75
$ws.Range("L46").Value = 84989.75
$ws.Range("N46").Value = -85451.75
$ws.Range("H132").Value = 4243.2197
$ws.Range("I132").Value = 3542.625
$ws.Range("J132").Value = 4910.452
$ws.Range("K132").Value = 10627.875
$ws.Range("L132").Value = 14731.356
$ws.Range("M132").Value = -8097.875
$ws.Range("N132").Value = -19791.356
$ws.Range("H134").Value = 84989.75
$ws.Range("J134").Value = 84989.75
$ws.Range("L134").Value = 254969.25
$ws.Range("N134").Value = -260039.25
$ws.Range("H136").Value = 3652.3542
$ws.Range("I136").Value = 3284.3462
$ws.Range("K136").Value = 9853.0386
$ws.Range("M136").Value = -7303.0386
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
